# Update the "Correspond Handoff Datetime" (E2) and
# "Correspond Handback DateTime" (H2) cells for the first data row
# on both the "zh-cn" and "de-de" worksheets, reflecting a re-run of
# the handback report generation.

$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-17 14:38:58"
$wsZh.Range("H2").Value = "2016-03-17 14:39:16"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-17 14:39:04"
$wsDe.Range("H2").Value = "2016-03-17 14:39:22"
